# "Se añade funcionalidad de comentarios"
# Fill in test case #10 (row 12, previously blank placeholder row), add two
# new test cases (#11 and #12, rows 13-14) documenting the new "comentarios"
# (comments) feature, and extend the sheet with a block of styled-but-empty
# rows (15-31) ready for future entries, matching the author's row template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting first --------------------------------------------------------
# Extend the styled table template (borders / centered wrap-text, plus the
# date format on column I) down through row 31, matching the look of the
# preceding rows, *before* putting values in so every new cell picks up the
# right style right away.
$ws.Range("H11:N11").Copy()
$ws.Range("H12:N31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 12 (test #10) -----------------------------------------------------
# Write text columns first, in the same order the author appears to have
# typed them (action -> expected result -> actual result -> obtained ->
# post-test action), so new shared-string entries land in the same order.
$ws.Range("J12").Value = "Intento agregar un comentario con la funcionalidad recientemente agregada"
$ws.Range("K12").Value = "Genera comentario y devuelve a Inicio"
$ws.Range("L12").Value = "ValueError"
$ws.Range("M12").Value = "NO"
$ws.Range("N12").Value = "Se debuguea en varias medidas el HTML, la URL, la view y se relaciona uno a uno la clase de comentarios para poder identificar cualquier objeto que se le brinda"
$ws.Range("H12").Value = 10
$ws.Range("I12").Value = 44995

# --- Row 13 (test #11) ------------------------------------------------------
$ws.Range("J13").Value = "Intento agregar un comentario a un objeto de cualquiera de las clases"
$ws.Range("K13").Value = "Genera comentario"
$ws.Range("L13").Value = "Genera comentario OK"
$ws.Range("M13").Value = "SI"
$ws.Range("N13").Value = "-"
$ws.Range("H13").Value = 11
$ws.Range("I13").Value = 44995

# --- Row 14 (test #12) ------------------------------------------------------
$ws.Range("K14").Value = "Eliminar comentario y quedarse en la página"
$ws.Range("L14").Value = "Elimina comentario OK"
$ws.Range("M14").Value = "SI"
$ws.Range("N14").Value = "-"
$ws.Range("H14").Value = 12
$ws.Range("I14").Value = 44995
$ws.Range("J14").Value = "Intento eliminar un comentario propio"

# --- Formatting ---------------------------------------------------------
# Row heights for the two wrapped-text rows.
$ws.Range("H12").RowHeight = 60
$ws.Range("H13").RowHeight = 30

# --- Conditional formatting --------------------------------------------------
# The "NO/SI" conditional formatting on column M now covers the newly added
# rows too (M9's separate rule, previously splitting the range in two, is no
# longer a gap once the table is filled in further down).
$fcs = $ws.Range("M3").FormatConditions
$newRange = $ws.Range("M3:M1048576")
$fcs.Item(1).ModifyAppliesToRange($newRange)
$fcs.Item(2).ModifyAppliesToRange($newRange)

# --- View state ---------------------------------------------------------
$ws.Range("J15").Select()
